$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 (I0) and J1 (IF), matching style of existing headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill I and J columns for data rows 2..39
# I column is a constant 1, except row 2 which is 8
# J column mirrors the existing H column (IP) value, except row 2 which is 8
for ($r = 2; $r -le 39; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    if ($r -eq 2) {
        $ws.Cells.Item($r, 9).Value = 8
        $ws.Cells.Item($r, 10).Value = 8
    } else {
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
